$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that currently sits at the very end of
#    the document (right after the final "};" run).
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2. Locate the paragraph containing the "it helps in cloning..." sentence so
#    we can edit the trailing whitespace / tab runs that follow it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*it helps in cloning*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# 3. Remove the single trailing space after "...updates on" (the run's text
#    loses its trailing space / "preserve" flag).
$find1 = $d.Range($pStart, $target.Range.End)
$find1.Find.Execute("updates on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterOn = $find1.End
$trailSpace = $d.Range($afterOn, $afterOn + 1)
$trailSpace.Text = ""

# 4. Delete the run of tab characters plus the first two of the three spaces
#    that used to precede "second obj ..." -- this collapses the seven
#    <w:tab/> runs down to nothing and leaves a single leading space on the
#    final run.
$find2 = $d.Range($afterOn, $target.Range.End)
$find2.Find.Execute("second obj", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$secondStart = $find2.Start
$gap = $d.Range($afterOn, $secondStart - 1)
$gap.Text = ""

# 5. Insert the "_GoBack" bookmark exactly between the two runs, i.e.
#    immediately after "...updates on" and before " second obj ...".
$bm = $d.Range($afterOn, $afterOn)
$d.Bookmarks.Add("_GoBack", $bm)
